# Apply value updates to the Kujata_Profits workbook (profit/price recalculations)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1012.9474
$ws.Range("J17").Value = 1012.9474
$ws.Range("L17").Value = 3038.8422
$ws.Range("N17").Value = -3374.8422
$ws.Range("H43").Value = 4283366.5
$ws.Range("I43").Value = 22800.2
$ws.Range("J43").Value = 6946220
$ws.Range("K43").Value = 22800.2
$ws.Range("L43").Value = 6946220
$ws.Range("M43").Value = -22731.2
$ws.Range("N43").Value = -6946358
$ws.Range("H58").Value = 1016.6
$ws.Range("J58").Value = 2504.25
$ws.Range("L58").Value = 7512.75
$ws.Range("N58").Value = -7812.75
$ws.Range("H87").Value = 39338.11
$ws.Range("J87").Value = 39338.11
$ws.Range("L87").Value = 39338.11
$ws.Range("N87").Value = -41834.11
$ws.Range("H90").Value = 39338.11
$ws.Range("J90").Value = 39338.11
$ws.Range("L90").Value = 118014.33
$ws.Range("N90").Value = -130494.33
$ws.Range("H132").Value = 5380972
$ws.Range("I132").Value = 7940854.5
$ws.Range("K132").Value = 23822563.5
$ws.Range("M132").Value = -23820033.5
$ws.Range("H135").Value = 23256290
$ws.Range("I135").Value = 214.43243
$ws.Range("K135").Value = 1929.89187
$ws.Range("M135").Value = 605.1081299999998
$ws.Range("H137").Value = 1183.7142
$ws.Range("I137").Value = 891.7778
$ws.Range("K137").Value = 2675.3334
$ws.Range("M137").Value = -125.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1702.9474
$ws.Range("I74").Value = 1221.8667
$ws.Range("K74").Value = 1221.8667
$ws.Range("M74").Value = -347.8667
$ws.Range("H77").Value = 1702.9474
$ws.Range("I77").Value = 1221.8667
$ws.Range("K77").Value = 6109.333500000001
$ws.Range("M77").Value = -1741.333500000001
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H132").Value = 1743.2703
$ws.Range("I132").Value = 1438.16
$ws.Range("J132").Value = 2378.9167
$ws.Range("K132").Value = 4314.48
$ws.Range("L132").Value = 7136.750100000001
$ws.Range("M132").Value = -1784.48
$ws.Range("N132").Value = -12196.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 13280
$ws.Range("J81").Value = 13280
$ws.Range("L81").Value = 13280
$ws.Range("N81").Value = -15402
$ws.Range("H84").Value = 13280
$ws.Range("J84").Value = 13280
$ws.Range("L84").Value = 39840
$ws.Range("N84").Value = -50448
$ws.Range("H86").Value = 2130.0356
$ws.Range("I86").Value = 2607.9
$ws.Range("J86").Value = 935.375
$ws.Range("K86").Value = 2607.9
$ws.Range("L86").Value = 935.375
$ws.Range("M86").Value = -1484.9
$ws.Range("N86").Value = -3181.375
$ws.Range("H88").Value = 30166.5
$ws.Range("J88").Value = 30166.5
$ws.Range("L88").Value = 30166.5
$ws.Range("N88").Value = -30978.5
$ws.Range("H89").Value = 2130.0356
$ws.Range("I89").Value = 2607.9
$ws.Range("J89").Value = 935.375
$ws.Range("K89").Value = 13039.5
$ws.Range("L89").Value = 4676.875
$ws.Range("M89").Value = -7423.5
$ws.Range("N89").Value = -15908.875
$ws.Range("H91").Value = 30166.5
$ws.Range("J91").Value = 30166.5
$ws.Range("L91").Value = 30166.5
$ws.Range("N91").Value = -32974.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1701.3823
$ws.Range("I31").Value = 1701.3823
$ws.Range("K31").Value = 1701.3823
$ws.Range("M31").Value = -1406.3823
$ws.Range("H34").Value = 1701.3823
$ws.Range("I34").Value = 1701.3823
$ws.Range("K34").Value = 1701.3823
$ws.Range("M34").Value = -1499.3823
$ws.Range("H132").Value = 1621.0238
$ws.Range("I132").Value = 1271.8485
$ws.Range("K132").Value = 3815.5455
$ws.Range("M132").Value = -1285.5455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 6282.0586
$ws.Range("I107").Value = 358.85715
$ws.Range("J107").Value = 10428.3
$ws.Range("K107").Value = 1076.57145
$ws.Range("L107").Value = 31284.9
$ws.Range("M107").Value = 843.4285500000001
$ws.Range("N107").Value = -35124.89999999999
$ws.Range("H125").Value = 5448.8335
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 6338.6
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 19015.8
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -28855.8
$ws.Range("H131").Value = 23812904
$ws.Range("J131").Value = 4150.8184
$ws.Range("L131").Value = 12452.4552
$ws.Range("N131").Value = -22532.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 641.4
$ws.Range("I107").Value = 977.4545000000001
$ws.Range("J107").Value = 377.35715
$ws.Range("K107").Value = 977.4545000000001
$ws.Range("L107").Value = 377.35715
$ws.Range("M107").Value = 942.5454999999999
$ws.Range("N107").Value = -4217.35715
$ws.Range("H121").Value = 9998
$ws.Range("J121").Value = 9998
$ws.Range("L121").Value = 9998
$ws.Range("N121").Value = -13492
$ws.Range("H132").Value = 2778.7778
$ws.Range("I132").Value = 2821.9285
$ws.Range("J132").Value = 2732.3076
$ws.Range("K132").Value = 8465.7855
$ws.Range("L132").Value = 8196.9228
$ws.Range("M132").Value = -5935.7855
$ws.Range("N132").Value = -13256.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2126.5334
$ws.Range("I82").Value = 2127.182
$ws.Range("J82").Value = 2124.75
$ws.Range("K82").Value = 2127.182
$ws.Range("L82").Value = 2124.75
$ws.Range("M82").Value = -1766.182
$ws.Range("N82").Value = -2846.75
$ws.Range("H85").Value = 2126.5334
$ws.Range("I85").Value = 2127.182
$ws.Range("J85").Value = 2124.75
$ws.Range("K85").Value = 2127.182
$ws.Range("L85").Value = 2124.75
$ws.Range("M85").Value = -879.1819999999998
$ws.Range("N85").Value = -4620.75
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = ""
$ws.Range("N119").Value = ""
$ws.Range("H141").Value = 70715
$ws.Range("J141").Value = 70715
$ws.Range("L141").Value = 70715
$ws.Range("N141").Value = -81075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20462
$ws.Range("H81").Value = 508.5
$ws.Range("I81").Value = 410.2
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 820.4
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = 240.6
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 508.5
$ws.Range("I84").Value = 410.2
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 4102
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = 1202
$ws.Range("N84").Value = -20608
$ws.Range("H107").Value = 715
$ws.Range("I107").Value = 653.3333
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1959.9999
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -39.99990000000003
$ws.Range("N107").Value = -6540
$ws.Range("H132").Value = 1714.0444
$ws.Range("I132").Value = 1229.7693
$ws.Range("J132").Value = 2376.7368
$ws.Range("K132").Value = 3689.3079
$ws.Range("L132").Value = 7130.2104
$ws.Range("M132").Value = -1159.3079
$ws.Range("N132").Value = -12190.2104
$ws.Range("H134").Value = 20000
$ws.Range("J134").Value = 20000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -65070
